$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 125: PCB part quantity increased from 1 to 20 and now tagged "PCB Assembly" ---
$ws.Range("G125").Value = 20
$ws.Range("H125").Value = 20
$ws.Range("J125").Value = 5.4
$ws.Range("K125").Value = "PCB Assembly"

# --- Row 128: PCB part quantity decreased from 4 to 1, no longer tagged "PCB Assembly" ---
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 1
$ws.Range("J128").Value = 1
$ws.Range("K128").Value = ""

# --- Row 130: PCB part quantity decreased from 20 to 4 (remains tagged "PCB Assembly") ---
$ws.Range("G130").Value = 4
$ws.Range("H130").Value = 4
$ws.Range("J130").Value = 1.64

# --- "Used in Sections" reference text reordered for rows 133-140 ---
$ws.Range("K133").Value = "Wheel Assembly, Corner Steering, Mechanical Integration, Differential Pivot, Rocker-Bogie, Body"
$ws.Range("K134").Value = "Differential Pivot, Wheel Assembly"
$ws.Range("K135").Value = "Differential Pivot, Wheel Assembly"
$ws.Range("K136").Value = "Differential Pivot, Wheel Assembly"
$ws.Range("K137").Value = "Differential Pivot, Head Assembly, Wheel Assembly, Rocker-Bogie"
$ws.Range("K138").Value = "Differential Pivot, Wheel Assembly"
$ws.Range("K139").Value = "Differential Pivot, Wheel Assembly, Rocker-Bogie, Mechanical Integration"
$ws.Range("K140").Value = "Differential Pivot, Head Assembly"
